$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Columns("BE:BE").Insert()
foreach ($n in $wb.Names) {
    Write-Output ($n.Name + " = " + $n.RefersTo)
}
